# "fix: new table db : tf_distribuicao_elegiveis; upload update new files ; tabela populada"
#
# The "Coordenação" abbreviation column (D) previously held the short
# acronym of the row's own Coordenação (e.g. "COFIS", "CMIF", ...), which
# duplicated information already present in column C and only added noise
# to the shared-string table. The column is repurposed to hold the
# abbreviation of the row's Diretoria (DIMAN / DISAT / DIBIO) instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Diretoria de Criação e Manejo de Unidades de Conservação - DIMAN
$ws.Range("D2:D5").Value = "DIMAN"

# Diretoria de Ações Socioambientais e Consolidação Territorial em
# Unidades de Conservação - DISAT
$ws.Range("D6").Value = "DISAT"
$ws.Range("D7").Value = "DISAT"
$ws.Range("D9").Value = "DISAT"

# Diretoria de Pesquisa, Avaliação e Monitoramento da Biodiversidade - DIBIO
$ws.Range("D8").Value = "DIBIO"
$ws.Range("D10").Value = "DIBIO"
$ws.Range("D11").Value = "DIBIO"
$ws.Range("D12").Value = "DIBIO"

# Widen column A (it now needs to comfortably fit the long Diretoria
# names) and mark it so Excel treats the width as content-fitted.
$ws.Columns("A").ColumnWidth = 90.16666666666667

# Move the active selection/view: the sheet no longer needs to be scrolled
# to column D on open, and the last user selection was cell B16.
$ws.Range("B16").Select() | Out-Null
